# edit.ps1 - applies the changes described by the target diff:
#  1. Slide 5, "Content Placeholder 1": merge the two runs of paragraph 2
#     ("The RAICHU API ... existing " + "source code.") into a single run.
#  2. Slide 5, "Title 2" ("THE IDEA"): apply the "Alien Encounters Solid"
#     typeface to the title run.
#  3. Slide 6, "Content Placeholder 1": split the run in paragraph 2
#     ("Whether is a lock ... than can be controlled ... RAICHU.") into
#     three runs, fixing the "than can" typo to "that needs to".

function Get-ShapeByName($slide, $name, $fallbackIndex) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $slide.Shapes.Item($fallbackIndex)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Change 1: slide 5 - merge the "source code." run into the previous run
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$contentShape5 = Get-ShapeByName $slide5 "Content Placeholder 1" 1
$body5 = $contentShape5.TextFrame.TextRange

$para2 = $body5.Paragraphs(2, 1)
# First overwrite with a placeholder so the engine treats the whole
# paragraph as genuinely changed (a no-op "identical text" assignment
# would otherwise leave the original two runs untouched/unmerged).
$para2.Text = "PLACEHOLDER_TEXT_RESET_0001"
$para2b = $body5.Paragraphs(2, 1)
$para2b.Text = "The RAICHU API will allow developers to easily implement the cloud services to new projects as well inject the capabilities into existing source code."

# ---------------------------------------------------------------------
# Change 2: slide 5 - title "THE IDEA" gets the Alien Encounters Solid font
# ---------------------------------------------------------------------
$titleShape5 = Get-ShapeByName $slide5 "Title 2" 2
$title5 = $titleShape5.TextFrame.TextRange
$title5.Font.Name = "Alien Encounters Solid"

# ---------------------------------------------------------------------
# Change 3: slide 6 - split "than can " into "that needs to " (3 runs)
# ---------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$contentShape6 = Get-ShapeByName $slide6 "Content Placeholder 1" 1
$body6 = $contentShape6.TextFrame.TextRange

$para2_6 = $body6.Paragraphs(2, 1)
$start6 = $para2_6.Start
$fullText6 = $para2_6.Text

$oldFragment = "than can "
$idx = $fullText6.IndexOf($oldFragment)
if ($idx -ge 0) {
    $seg = $body6.Characters($start6 + $idx, $oldFragment.Length)
    $seg.Text = "that needs to "
}
